$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# Column A holds the "card" number as text (matching the other rows/sheets,
# e.g. A2 already = "9"). Use a leading apostrophe so Excel stores it as
# text "9" instead of re-inferring it as a number.
for ($r = 3; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "'9"
}
